$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 'Supplemental Digital Content is available in the text.
'
$ws.Cells.Item(2, 5).Value = '[Sapna R.%Kudchadkar%NULL%0, Christopher L.%Carroll%NULL%2]'
$ws.Cells.Item(2, 9).Value = ''
$ws.Cells.Item(2, 10).Value = 'Lippincott Williams &amp; Wilkins'

# Row 3
$ws.Cells.Item(3, 5).Value = '[Sohaib R%Rufai%sohaibrufai@gmail.com%0, Catey%Bunce%Catey.bunce@kcl.ac.uk%2, Catey%Bunce%Catey.bunce@kcl.ac.uk%0]'
$ws.Cells.Item(3, 9).Value = ''
$ws.Cells.Item(3, 10).Value = 'Oxford University Press'

# Row 4
$ws.Cells.Item(4, 5).Value = '[Travis%Sanchez%NULL%0, Kamal%Al Nasr%NULL%2, Kamal%Al Nasr%NULL%0, Ross%Gore%NULL%1, Abdullah%Wahbeh%NULL%2, Abdullah%Wahbeh%NULL%0, Tareq%Nasralah%t.nasralah@northeastern.edu%2, Tareq%Nasralah%t.nasralah@northeastern.edu%0, Mohammad%Al-Ramahi%NULL%2, Mohammad%Al-Ramahi%NULL%0, Omar%El-Gayar%NULL%2, Omar%El-Gayar%NULL%0]'
$ws.Cells.Item(4, 9).Value = ''
$ws.Cells.Item(4, 10).Value = 'JMIR Publications'

# Row 5
$ws.Cells.Item(5, 5).Value = '[Konstantin S.%Sharov%const.sharov@mail.ru%0]'
$ws.Cells.Item(5, 9).Value = ''
$ws.Cells.Item(5, 10).Value = 'John Wiley and Sons Inc.'

# Row 6
$ws.Cells.Item(6, 3).Value = 'Unknown Title'
$ws.Cells.Item(6, 6).Value = 'not found'
$ws.Cells.Item(6, 7).Value = 'N/A'
$ws.Cells.Item(6, 8).Value = '1970-01-01'
$ws.Cells.Item(6, 10).Value = ''

# Row 7
$ws.Cells.Item(7, 3).Value = 'Unknown Title'
$ws.Cells.Item(7, 4).Value = 'Unknown Abstract'
$ws.Cells.Item(7, 5).Value = '[]'
$ws.Cells.Item(7, 6).Value = 'not found'
$ws.Cells.Item(7, 7).Value = 'N/A'
$ws.Cells.Item(7, 8).Value = '1970-01-01'
$ws.Cells.Item(7, 10).Value = ''

# Row 8
$ws.Cells.Item(8, 3).Value = 'Unknown Title'
$ws.Cells.Item(8, 4).Value = 'Unknown Abstract'
$ws.Cells.Item(8, 5).Value = '[]'
$ws.Cells.Item(8, 6).Value = 'not found'
$ws.Cells.Item(8, 7).Value = 'N/A'
$ws.Cells.Item(8, 8).Value = '1970-01-01'
$ws.Cells.Item(8, 10).Value = ''

# Row 9
$ws.Cells.Item(9, 5).Value = '[Vittorio%Gebbia%NULL%0, Dario%Piazza%NULL%1, Maria Rosaria%Valerio%NULL%1, Nicolò%Borsellino%NULL%1, Alberto%Firenze%NULL%1]'
$ws.Cells.Item(9, 9).Value = ''
$ws.Cells.Item(9, 10).Value = 'American Society of Clinical Oncology'

# Row 10
$ws.Cells.Item(10, 3).Value = 'Unknown Title'
$ws.Cells.Item(10, 5).Value = '[]'
$ws.Cells.Item(10, 6).Value = 'not found'
$ws.Cells.Item(10, 7).Value = 'N/A'
$ws.Cells.Item(10, 8).Value = '1970-01-01'
$ws.Cells.Item(10, 10).Value = ''

# Row 11
$ws.Cells.Item(11, 5).Value = '[Maria Renee%Jimenez‐Sotomayor%NULL%0, Carolina%Gomez‐Moreno%NULL%1, Enrique%Soto‐Perez‐de‐Celis%enrique.sotop@incmnsz.mx%1]'
$ws.Cells.Item(11, 9).Value = ''
$ws.Cells.Item(11, 10).Value = 'John Wiley &amp; Sons, Inc.'

# Row 12
$ws.Cells.Item(12, 4).Value = 'Background
id="Par1">Internet analytics are increasingly being integrated into public health regulation.

 One specific application is to monitor compliance of website and social media activity with respect to jurisdictional regulations.

 These data may then identify breaches of compliance and inform disciplinary actions.

 Our study aimed to evaluate the novel use of internet analytics by a Canadian chiropractic regulator to determine their registrants compliance with three regulations related to specific health conditions, pregnancy conditions and most recently, claims of improved immunity during the COVID-19 crisis.


Methods
id="Par2">A customized internet search tool (Market Review Tool, MRT) was used by the College of Chiropractors of British Columbia (CCBC), Canada to audit registrants websites and social media activity.

 The audits extracted words whose use within specific contexts is not permitted under CCBC guidelines.

 The MRT was first used in October of 2018 to identify words related to specific health conditions.

 The MRT was again used in December 2019 for words related to pregnancy and most recently in March 2020 for words related to COVID-19. In these three MRT applications, potential cases of word misuse were evaluated by the regulator who then notified the practitioner to comply with existing regulations by a specific date.

 The MRT was then used on that date to determine compliance.

 Those found to be non-compliant were referred to the regulator’s inquiry committee.

 We mapped this process and reported the outcomes with permission of the regulator.


Results
id="Par3">In September 2018, 250 inappropriate mentions of specific health conditions were detected from approximately 1250 registrants with 2 failing to comply.

 The second scan for pregnancy related terms of approximately1350 practitioners revealed 83 inappropriate mentions.

 Following notification, all 83 cases were compliant within the specified timeframe.

 Regarding COVID-19 related words, 97 inappropriate mentions of the word “immune” were detected from 1350 registrants with 7 cases of non-compliance.


Conclusion
id="Par4">Internet analytics are an effective way for regulators to monitor internet activity to protect the public from misleading statements.

 The processes described were effective at bringing about rapid practitioner compliance.

 Given the increasing volume of internet activity by healthcare professionals, internet analytics are an important addition for health care regulators to protect the public they serve.


'
$ws.Cells.Item(12, 5).Value = '[Greg%Kawchuk%greg.kawchuk@ualberta.ca%0, Jan%Hartvigsen%jhartvigsen@health.sdu.dk%2, Jan%Hartvigsen%jhartvigsen@health.sdu.dk%0, Stan%Innes%S.Innes@murdoch.edu.au%1, J. Keith%Simpson%k.simpson@murdoch.edu.au%1, Brian%Gushaty%bgushaty@gushaty.com%1]'
$ws.Cells.Item(12, 9).Value = ''
$ws.Cells.Item(12, 10).Value = 'BioMed Central'

# Row 13
$ws.Cells.Item(13, 5).Value = '[Gunther%Eysenbach%NULL%0, Jon-Patrick%Allem%NULL%2, Jon-Patrick%Allem%NULL%0, Richard%Zowalla%NULL%1, Wasim%Ahmed%Wasim.Ahmed@Newcastle.ac.uk%2, Wasim%Ahmed%Wasim.Ahmed@Newcastle.ac.uk%0, Josep%Vidal-Alaball%NULL%2, Josep%Vidal-Alaball%NULL%0, Joseph%Downing%NULL%2, Joseph%Downing%NULL%0, Francesc%López Seguí%NULL%2, Francesc%López Seguí%NULL%0]'
$ws.Cells.Item(13, 9).Value = ''
$ws.Cells.Item(13, 10).Value = 'JMIR Publications'

# Row 14
$ws.Cells.Item(14, 3).Value = 'Unknown Title'
$ws.Cells.Item(14, 4).Value = 'Unknown Abstract'
$ws.Cells.Item(14, 6).Value = 'not found'
$ws.Cells.Item(14, 7).Value = 'N/A'
$ws.Cells.Item(14, 8).Value = '1970-01-01'
$ws.Cells.Item(14, 10).Value = ''

# Row 15
$ws.Cells.Item(15, 3).Value = 'Unknown Title'
$ws.Cells.Item(15, 4).Value = 'Unknown Abstract'
$ws.Cells.Item(15, 5).Value = '[]'
$ws.Cells.Item(15, 6).Value = 'not found'
$ws.Cells.Item(15, 7).Value = 'N/A'
$ws.Cells.Item(15, 8).Value = '1970-01-01'
$ws.Cells.Item(15, 10).Value = ''

# Row 16
$ws.Cells.Item(16, 4).Value = 'id="Par1">The aim of this study was to elicit the views of medical faculty students regarding the COVID-19 pandemic.
 This descriptive study was performed with Ondokuz Mayıs University Medical Faculty students on 24–27 March, 2020. The Medical Faculty currently has 2051 students.
 A questionnaire was used as a data collection tool.
 For that purpose, the authors designed a questionnaire specifically for this research via the “Google Forms” web.
 This consisted of 40 open- and close-ended questions.
 The questionnaire was completed by 1375 (67.1%) students.
 Accordingly, 52.4% of medical students reported feeling mentally unwell.
 Although 50.8% of medical students reported generally/usually obtaining information about COVID-19 through the social media, 82.0% did not trust information/messages arriving through the social media and WhatsApp.
 We found that 86.7% of students regarded frequent hand washing as the most important means of protection against COVID-19, and 19.3% of students did not regard the COVID-19 pandemic as a severe public health problem for Turkey at that moment.
 In addition, 61.6% of students stated that a suppression strategy involving tight restrictions need to be applied to being the pandemic under control in Turkey.
 Use can be made of medical students in the transmission of accurate information during the COVID-19 pandemic.
 Students can be excellent activists on these subjects in countries in which medical education is suspended.
 Measures therefore need to be taken concerning the transmission of up to date and accurate information to medical students.
'
$ws.Cells.Item(16, 5).Value = '[Servet%Aker%servetaker@gmail.com%0, Özlem%Mıdık%NULL%2, Özlem%Mıdık%NULL%0]'
$ws.Cells.Item(16, 9).Value = ''
$ws.Cells.Item(16, 10).Value = 'Springer US'

# Row 17
$ws.Cells.Item(17, 3).Value = 'Unknown Title'
$ws.Cells.Item(17, 4).Value = 'Unknown Abstract'
$ws.Cells.Item(17, 5).Value = '[]'
$ws.Cells.Item(17, 6).Value = 'not found'
$ws.Cells.Item(17, 7).Value = 'N/A'
$ws.Cells.Item(17, 8).Value = '1970-01-01'
$ws.Cells.Item(17, 10).Value = ''

# Row 18
$ws.Cells.Item(18, 3).Value = 'Unknown Title'
$ws.Cells.Item(18, 4).Value = 'Unknown Abstract'
$ws.Cells.Item(18, 5).Value = '[]'
$ws.Cells.Item(18, 6).Value = 'not found'
$ws.Cells.Item(18, 7).Value = 'N/A'
$ws.Cells.Item(18, 9).Value = ''

# Row 19
$ws.Cells.Item(19, 5).Value = '[Bahar%Yuksel%baharyl86@gmail.com%0, Kubra%Cakmak%NULL%1]'
$ws.Cells.Item(19, 9).Value = ''
$ws.Cells.Item(19, 10).Value = 'John Wiley and Sons Inc.'

# Row 20
$ws.Cells.Item(20, 5).Value = '[Gunther%Eysenbach%NULL%0, Nazakat%Hamassed%NULL%2, Nazakat%Hamassed%NULL%0, Hardawan%Kakashekh%NULL%1, Muhammad%Saud%NULL%1, Mohammad Amin%Bahrami%NULL%1, Araz Ramazan%Ahmad%araz.ahmad85@uor.edu.krd%2, Araz Ramazan%Ahmad%araz.ahmad85@uor.edu.krd%0, Hersh Rasool%Murad%NULL%2, Hersh Rasool%Murad%NULL%0]'
$ws.Cells.Item(20, 9).Value = ''
$ws.Cells.Item(20, 10).Value = 'JMIR Publications'

# Row 21
$ws.Cells.Item(21, 5).Value = '[Alexander%Muacevic%NULL%0, John R%Adler%NULL%0, Ramez%Kouzy%NULL%2, Ramez%Kouzy%NULL%0, Joseph%Abi Jaoude%NULL%1, Afif%Kraitem%NULL%1, Molly B%El Alam%NULL%1, Basil%Karam%NULL%1, Elio%Adib%NULL%1, Jabra%Zarka%NULL%1, Cindy%Traboulsi%NULL%1, Elie W%Akl%NULL%1, Khalil%Baddour%NULL%1]'
$ws.Cells.Item(21, 9).Value = ''
$ws.Cells.Item(21, 10).Value = 'Cureus'
